# Weekly update: insert a new price record at row 379 (date 2023-06-16),
# shifting all subsequent rows (379-495) down by one to (380-496).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 379; Excel shifts rows 379..495 down to 380..496
# and inherits the formatting from the row above (so the D column keeps its date style).
$ws.Rows.Item(379).Insert()

# Populate the newly inserted row 379 with the new weekly price record.
$ws.Range("A379").Value = 5
$ws.Range("B379").Value = "Macroferia Regional de Talca"
$ws.Range("C379").Value = "Maule"
$ws.Range("D379").Value = 45093
$ws.Range("E379").Value = 7
$ws.Range("F379").Value = 100112003
$ws.Range("G379").Value = "Ajo"
$ws.Range("H379").Value = "Chino"
$ws.Range("I379").Value = "Primera"
$ws.Range("J379").Value = 300
$ws.Range("K379").Value = 18000
$ws.Range("L379").Value = 18000
$ws.Range("M379").Value = 18000
$ws.Range("N379").Value = "`$/malla 10 kilos"
$ws.Range("O379").Value = "China"
$ws.Range("P379").Value = 1800
$ws.Range("Q379").Value = 10
$ws.Range("R379").Value = "Hortaliza"
